$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the new row the same formatting as the row above it (the last line
# of the existing "Sparkfun" vendor block) before filling in its values.
[void]$ws.Range("C11").Copy()
[void]$ws.Range("C12").PasteSpecial(-4122)   # xlPasteFormats
[void]$ws.Range("D11").Copy()
[void]$ws.Range("D12").PasteSpecial(-4122)
[void]$ws.Range("G11").Copy()
[void]$ws.Range("G12").PasteSpecial(-4122)

# New row 12: LiPo Fuel Gauge
$ws.Range("A12").Value = "LiPo Fuel Gauge"
$ws.Range("B12").Value = "https://www.sparkfun.com/products/10617"
$ws.Range("D12").Value = 9.95
$ws.Range("E12").Value = 2
$ws.Range("G12").Formula = "=D12*E12+F12"

# The "Sparkfun" vendor column was merged C9:C11; extend it to cover the
# newly added row.
[void]$ws.Range("C9:C12").Merge()

# Move the active selection down to the new row, like the author did.
[void]$ws.Range("H12").Select()
